# The "Version Control" table (the document's only table) lists, per
# revision row, a "responsible" person and a "reviewer" person, each stored
# as "<Name><space>(<ROLE>)" split across two runs: one run holding the name
# (optionally with a trailing space) and a second run holding "(ROLE)"
# (optionally with a leading space). This change reassigns some of those
# names/roles for a few revision rows:
#
#   Row 6 (v2.8.1) col 5 (reviewer)    : กิตติพศ (SP) -> วิรัตน์ (TL)
#   Row 7 (v1.4.1) col 4 (responsible) : ณัฐดนัย (DM) -> กิตติพศ (SP)
#   Row 7 (v1.4.1) col 5 (reviewer)    : กิตติพศ (SP) -> วิรัตน์ (TL)
#   Row 8 (v1.2.2) col 4 (responsible) : ณัฐดนัย (DM) -> วิรัตน์ (TL)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-NameRoleCell($row, $col, $newName, $newRole) {
    $cell = $t.Cell($row, $col)
    $rg = $cell.Range

    # First "word" of the cell is the name run (e.g. "กิตติพศ " or
    # "ณัฐดนัย"); everything between there and the cell's end-of-cell mark
    # (the range's final "word") is the role run (e.g. "(SP)" / " (DM)").
    $nameStart = $rg.Words.Item(1).Start
    $nameEnd = $rg.Words.Item(1).End
    $roleStart = $nameEnd
    $roleEnd = $rg.Words.Item($rg.Words.Count - 1).End

    # Write the role run before the name run: both edits target the same
    # cell, and writing back-to-front means the (not-yet-written) name run's
    # offsets are never invalidated by the role-run edit, and the two runs
    # stay distinct instead of merging into one.
    $d.Range($roleStart, $roleEnd).Text = $newRole
    $d.Range($nameStart, $nameEnd).Text = $newName
}

# Apply bottom-to-top (row 8, then 7, then 6): every row sits entirely
# before the next one character-wise, so finishing a lower row first keeps
# the as-yet-unprocessed rows' offsets untouched.
Set-NameRoleCell 8 4 "วิรัตน์" " (TL)"
Set-NameRoleCell 7 5 "วิรัตน์" " (TL)"
Set-NameRoleCell 7 4 "กิตติพศ " "(SP)"
Set-NameRoleCell 6 5 "วิรัตน์" " (TL)"
